$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row 5 (TestCaseID C02) - written first to match original authoring order
$ws.Range("A5").Value = "C02"

# New header cells for columns K:T on row 1, interleaved with row 5 values
# in the same order the original author typed them
$ws.Range("K1").Value = "GameTheme"
$ws.Range("K5").Value = "FAST"
$ws.Range("L1").Value = "GameTitle"
$ws.Range("M1").Value = "Topper"
$ws.Range("O1").Value = "Denom"
$ws.Range("P1").Value = "VAR"
$ws.Range("Q1").Value = "CandleColor"
$ws.Range("R1").Value = "TicketValidation"
$ws.Range("S1").Value = "AuditSwitch"
$ws.Range("N1").Value = "ButtonPanel"
$ws.Range("T1").Value = "UnitDisc"

$ws.Range("L5").Value = "TIGER WEALTH - (GT-TIGER WEALTH)"
$ws.Range("M5").Value = "TOPPER 19.5 INCH LCD HELIX XT LATAM - (TPP000127-LATAM|01)"
$ws.Range("N5").Value = "BTN PANEL 13 LCD HELIX SLANT AND UPRIGHT REFURB - (TBTN000079|01)"
$ws.Range("O5").Value = ".01 - (DENOM_3)"
$ws.Range("P5").Value = "0 - (VAR_0)"
$ws.Range("Q5").Value = "FILM CANDLE BLANK PLATE CHAMPAGNE - (CAN000015|01)"
$ws.Range("R5").Value = "NONE - USE OPTION WHEN THERE IS NO TKT - (TKV000011|01)"
$ws.Range("S5").Value = "2341 Audit Switch - IGT Flat Key - (AUDSW_3)"
$ws.Range("T5").Value = 100

# Match the saved selection/view state
[void]$ws.Range("S5").Select()
$excel.ActiveWindow.ScrollColumn = 6
